# log_201608040212_zouguo.xlsx update
# - fix day-3/4/5 time typos (4:31/4:32/4:33 -> 4:30)
# - day-6: change time window, record execution + follow-up notes
# - day-7: change time window, record work done / status / notes
# - add day-8 placeholder rows (22-28)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Day 3 (row 17): correct the typo'd end time
$ws.Range("A17").Value = "第三天1:00-4:30"

# Day 4 (row 18): correct the typo'd end time
$ws.Range("A18").Value = "第四天1:00-4:30"

# Day 5 (row 19): correct the typo'd end time
$ws.Range("A19").Value = "第五天1:00-4:30"

# Day 6 (row 20): new time window + follow-up plan note
$ws.Range("A20").Value = "  第六天10:00-14:30"
$ws.Range("D20").Value = "编写需求规格说明书"

# Day 7 (row 21): new time window + what was done / status / notes
$ws.Range("A21").Value = " 第七天9:00-17:00"
$ws.Range("B21").Value = "完善dao service 测试 以及完善代码注释"
$ws.Range("C21").Value = "完成"
$ws.Range("D21").Value = "编写需求规格说明书"

# Day 8 placeholder rows
$ws.Range("A22").Value = "第八天"
$ws.Range("A23").Value = "第八天"
$ws.Range("A24").Value = "第八天"
$ws.Range("A25").Value = "第八天"
$ws.Range("A26").Value = "第八天"
$ws.Range("A27").Value = "第八天"
$ws.Range("A28").Value = "第八天"

$ws.Range("B22").Select()
